$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header labels for the last "Prob iMOE<1" block: Min/Max -> 5th/95th
$ws.Range("Z2").Value = "5th"
$ws.Range("AA2").Value = "95th"

# Row 3 (Truong / Bastos-Moreira et al. 2023) updated Prob iMOE<1 stats
$ws.Range("W3").Value = 3.64
$ws.Range("X3").Value = 0.04
$ws.Range("Y3").Value = 18.420000000000002
$ws.Range("Z3").Value = 0
$ws.Range("AA3").Value = 46.64

# Row 4 (Plasma Groningen) updated Prob iMOE<1 stats
$ws.Range("W4").Value = 1.59
$ws.Range("X4").Value = 0.6
$ws.Range("Y4").Value = 5.39
$ws.Range("Z4").Value = 0
$ws.Range("AA4").Value = 22.15

# Row 5 (Plasma Burkina toddlers) updated Prob iMOE<1 stats
$ws.Range("W5").Value = 0
$ws.Range("X5").Value = 0
$ws.Range("Y5").Value = 0.66
$ws.Range("Z5").Value = 0
$ws.Range("AA5").Value = 12.48

# Row 6 (De ruyck urine) updated Prob iMOE<1 stats
$ws.Range("AA6").Value = 0.09

# Row 7 (Asam urine) updated Prob iMOE<1 stats
$ws.Range("AA7").Value = 0.13

# Update selection to A5
$ws.Range("A5").Select()
